# Apply updated TPM-derived values to the Adam9-Itgav LR-pairs sheet.
# The sheet contains a 5x5 (sending cluster x target cluster) table of
# ligand/receptor expression and specificity metrics. This script writes
# the refreshed values (ligand/receptor average & total expression,
# derived specificities, and their products -- edge weights/specificities)
# into the corresponding cells, row by row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 9.462749333333333
$ws.Range("H2").Value = 28.388248
$ws.Range("I2").Value = 0.07254428564686972
$ws.Range("J2").Value = 0.07439525120506714
$ws.Range("M2").Value = 13.89934866666667
$ws.Range("N2").Value = 41.69804600000001
$ws.Range("O2").Value = 0.04853507553134179
$ws.Range("P2").Value = 0.04999273878390351
$ws.Range("Q2").Value = 131.5260523292676
$ws.Range("R2").Value = 1183.734470963408
$ws.Range("S2").Value = 0.003520942383238056
$ws.Range("T2").Value = 0.003719222360257804
$ws.Range("G3").Value = 9.462749333333333
$ws.Range("H3").Value = 28.388248
$ws.Range("I3").Value = 0.07254428564686972
$ws.Range("J3").Value = 0.07439525120506714
$ws.Range("O3").Value = 0.245697991654417
$ws.Range("P3").Value = 0.253077086664408
$ws.Range("Q3").Value = 665.8212963255155
$ws.Range("R3").Value = 5992.39166692964
$ws.Range("S3").Value = 0.01782398528944024
$ws.Range("T3").Value = 0.01882773343664518
$ws.Range("G4").Value = 9.462749333333333
$ws.Range("H4").Value = 28.388248
$ws.Range("I4").Value = 0.07254428564686972
$ws.Range("J4").Value = 0.07439525120506714
$ws.Range("M4").Value = 82.007665
$ws.Range("N4").Value = 246.022995
$ws.Range("O4").Value = 0.2863622109480123
$ws.Range("P4").Value = 0.2949625822722868
$ws.Range("Q4").Value = 776.0179773069733
$ws.Range("R4").Value = 6984.16179576276
$ws.Range("S4").Value = 0.02077394202948177
$ws.Range("T4").Value = 0.02194381540424206
$ws.Range("G5").Value = 9.462749333333333
$ws.Range("H5").Value = 28.388248
$ws.Range("I5").Value = 0.07254428564686972
$ws.Range("J5").Value = 0.07439525120506714
$ws.Range("M5").Value = 25.0501465
$ws.Range("N5").Value = 50.100293
$ws.Range("O5").Value = 0.0874724982879541
$ws.Range("P5").Value = 0.06006638442832619
$ws.Range("Q5").Value = 237.0432570927773
$ws.Range("R5").Value = 1422.259542556664
$ws.Range("S5").Value = 0.006345629902046665
$ws.Range("T5").Value = 0.00446865375852546
$ws.Range("G6").Value = 9.462749333333333
$ws.Range("H6").Value = 28.388248
$ws.Range("I6").Value = 0.07254428564686972
$ws.Range("J6").Value = 0.07439525120506714
$ws.Range("M6").Value = 95.05788666666668
$ws.Range("N6").Value = 285.17366
$ws.Range("O6").Value = 0.3319322235782747
$ws.Range("P6").Value = 0.3419012078510756
$ws.Range("Q6").Value = 899.5089536830756
$ws.Range("R6").Value = 8095.58058314768
$ws.Range("S6").Value = 0.02407978604266298
$ws.Range("T6").Value = 0.02543582624539665
$ws.Range("I7").Value = 0.3231336970688258
$ws.Range("J7").Value = 0.3313784449305509
$ws.Range("M7").Value = 13.89934866666667
$ws.Range("N7").Value = 41.69804600000001
$ws.Range("O7").Value = 0.04853507553134179
$ws.Range("P7").Value = 0.04999273878390351
$ws.Range("Q7").Value = 585.8559247093223
$ws.Range("R7").Value = 5272.703322383901
$ws.Range("S7").Value = 0.01568331839395718
$ws.Range("T7").Value = 0.01656651603602918
$ws.Range("I8").Value = 0.3231336970688258
$ws.Range("J8").Value = 0.3313784449305509
$ws.Range("O8").Value = 0.245697991654417
$ws.Range("P8").Value = 0.253077086664408
$ws.Range("S8").Value = 0.07939330040567728
$ws.Range("T8").Value = 0.08386429142640577
$ws.Range("I9").Value = 0.3231336970688258
$ws.Range("J9").Value = 0.3313784449305509
$ws.Range("M9").Value = 82.007665
$ws.Range("N9").Value = 246.022995
$ws.Range("O9").Value = 0.2863622109480123
$ws.Range("P9").Value = 0.2949625822722868
$ws.Range("Q9").Value = 3456.613512189083
$ws.Range("R9").Value = 31109.52160970175
$ws.Range("S9").Value = 0.09253327992443422
$ws.Range("T9").Value = 0.09774424182609005
$ws.Range("I10").Value = 0.3231336970688258
$ws.Range("J10").Value = 0.3313784449305509
$ws.Range("M10").Value = 25.0501465
$ws.Range("N10").Value = 50.100293
$ws.Range("O10").Value = 0.0874724982879541
$ws.Range("P10").Value = 0.06006638442832619
$ws.Range("Q10").Value = 1055.860752457908
$ws.Range("R10").Value = 6335.16451474745
$ws.Range("S10").Value = 0.02826531176363315
$ws.Range("T10").Value = 0.01990470506445939
$ws.Range("I11").Value = 0.3231336970688258
$ws.Range("J11").Value = 0.3313784449305509
$ws.Range("M11").Value = 95.05788666666668
$ws.Range("N11").Value = 285.17366
$ws.Range("O11").Value = 0.3319322235782747
$ws.Range("P11").Value = 0.3419012078510756
$ws.Range("Q11").Value = 4006.678832913223
$ws.Range("R11").Value = 36060.109496219
$ws.Range("S11").Value = 0.107258486581124
$ws.Range("T11").Value = 0.1132986905775665
$ws.Range("G12").Value = 32.300192
$ws.Range("H12").Value = 96.900576
$ws.Range("I12").Value = 0.2476229975407503
$ws.Range("J12").Value = 0.2539410918713864
$ws.Range("M12").Value = 13.89934866666667
$ws.Range("N12").Value = 41.69804600000001
$ws.Range("O12").Value = 0.04853507553134179
$ws.Range("P12").Value = 0.04999273878390351
$ws.Range("Q12").Value = 448.9516306082774
$ws.Range("R12").Value = 4040.564675474497
$ws.Range("S12").Value = 0.01201840088893758
$ws.Range("T12").Value = 0.01269521067242546
$ws.Range("G13").Value = 32.300192
$ws.Range("H13").Value = 96.900576
$ws.Range("I13").Value = 0.2476229975407503
$ws.Range("J13").Value = 0.2539410918713864
$ws.Range("O13").Value = 0.245697991654417
$ws.Range("P13").Value = 0.253077086664408
$ws.Range("Q13").Value = 2272.717468404854
$ws.Range("R13").Value = 20454.45721564368
$ws.Range("S13").Value = 0.06084047318320899
$ws.Range("T13").Value = 0.06426667171518925
$ws.Range("G14").Value = 32.300192
$ws.Range("H14").Value = 96.900576
$ws.Range("I14").Value = 0.2476229975407503
$ws.Range("J14").Value = 0.2539410918713864
$ws.Range("M14").Value = 82.007665
$ws.Range("N14").Value = 246.022995
$ws.Range("O14").Value = 0.2863622109480123
$ws.Range("P14").Value = 0.2949625822722868
$ws.Range("Q14").Value = 2648.86332497168
$ws.Range("R14").Value = 23839.76992474512
$ws.Range("S14").Value = 0.07090986905734348
$ws.Range("T14").Value = 0.07490312020342815
$ws.Range("G15").Value = 32.300192
$ws.Range("H15").Value = 96.900576
$ws.Range("I15").Value = 0.2476229975407503
$ws.Range("J15").Value = 0.2539410918713864
$ws.Range("M15").Value = 25.0501465
$ws.Range("N15").Value = 50.100293
$ws.Range("O15").Value = 0.0874724982879541
$ws.Range("P15").Value = 0.06006638442832619
$ws.Range("Q15").Value = 809.1245415781281
$ws.Range("R15").Value = 4854.747249468768
$ws.Range("S15").Value = 0.02166020222844134
$ws.Range("T15").Value = 0.01525332324649559
$ws.Range("G16").Value = 32.300192
$ws.Range("H16").Value = 96.900576
$ws.Range("I16").Value = 0.2476229975407503
$ws.Range("J16").Value = 0.2539410918713864
$ws.Range("M16").Value = 95.05788666666668
$ws.Range("N16").Value = 285.17366
$ws.Range("O16").Value = 0.3319322235782747
$ws.Range("P16").Value = 0.3419012078510756
$ws.Range("Q16").Value = 3070.387990447574
$ws.Range("R16").Value = 27633.49191402816
$ws.Range("S16").Value = 0.08219405218281889
$ws.Range("T16").Value = 0.08682276603384798
$ws.Range("G17").Value = 9.736177999999999
$ws.Range("H17").Value = 19.472356
$ws.Range("I17").Value = 0.07464047213559308
$ws.Range("J17").Value = 0.0510299478916239
$ws.Range("M17").Value = 13.89934866666667
$ws.Range("N17").Value = 41.69804600000001
$ws.Range("O17").Value = 0.04853507553134179
$ws.Range("P17").Value = 0.04999273878390351
$ws.Range("Q17").Value = 135.3265327027293
$ws.Range("R17").Value = 811.959196216376
$ws.Range("S17").Value = 0.003622680952796022
$ws.Range("T17").Value = 0.002551126855102161
$ws.Range("G18").Value = 9.736177999999999
$ws.Range("H18").Value = 19.472356
$ws.Range("I18").Value = 0.07464047213559308
$ws.Range("J18").Value = 0.0510299478916239
$ws.Range("O18").Value = 0.245697991654417
$ws.Range("P18").Value = 0.253077086664408
$ws.Range("Q18").Value = 685.0603803252633
$ws.Range("R18").Value = 4110.362281951579
$ws.Range("S18").Value = 0.0183390140998527
$ws.Range("T18").Value = 0.01291451054504873
$ws.Range("G19").Value = 9.736177999999999
$ws.Range("H19").Value = 19.472356
$ws.Range("I19").Value = 0.07464047213559308
$ws.Range("J19").Value = 0.0510299478916239
$ws.Range("M19").Value = 82.007665
$ws.Range("N19").Value = 246.022995
$ws.Range("O19").Value = 0.2863622109480123
$ws.Range("P19").Value = 0.2949625822722868
$ws.Range("Q19").Value = 798.4412238043699
$ws.Range("R19").Value = 4790.64734282622
$ws.Range("S19").Value = 0.02137421062695194
$ws.Range("T19").Value = 0.01505192520333362
$ws.Range("G20").Value = 9.736177999999999
$ws.Range("H20").Value = 19.472356
$ws.Range("I20").Value = 0.07464047213559308
$ws.Range("J20").Value = 0.0510299478916239
$ws.Range("M20").Value = 25.0501465
$ws.Range("N20").Value = 50.100293
$ws.Range("O20").Value = 0.0874724982879541
$ws.Range("P20").Value = 0.06006638442832619
$ws.Range("Q20").Value = 243.892685250077
$ws.Range("R20").Value = 975.5707410003079
$ws.Range("S20").Value = 0.006528988571092752
$ws.Range("T20").Value = 0.003065184467415734
$ws.Range("G21").Value = 9.736177999999999
$ws.Range("H21").Value = 19.472356
$ws.Range("I21").Value = 0.07464047213559308
$ws.Range("J21").Value = 0.0510299478916239
$ws.Range("M21").Value = 95.05788666666668
$ws.Range("N21").Value = 285.17366
$ws.Range("O21").Value = 0.3319322235782747
$ws.Range("P21").Value = 0.3419012078510756
$ws.Range("Q21").Value = 925.5005048904933
$ws.Range("R21").Value = 5553.00302934296
$ws.Range("S21").Value = 0.02477557788489966
$ws.Range("T21").Value = 0.01744720082072366
$ws.Range("G22").Value = 36.79199966666667
$ws.Range("H22").Value = 110.375999
$ws.Range("I22").Value = 0.2820585476079611
$ws.Range("J22").Value = 0.2892552641013719
$ws.Range("M22").Value = 13.89934866666667
$ws.Range("N22").Value = 41.69804600000001
$ws.Range("O22").Value = 0.04853507553134179
$ws.Range("P22").Value = 0.04999273878390351
$ws.Range("Q22").Value = 511.3848315108839
$ws.Range("R22").Value = 4602.463483597955
$ws.Range("S22").Value = 0.01368973291241296
$ws.Range("T22").Value = 0.01446066286008891
$ws.Range("G23").Value = 36.79199966666667
$ws.Range("H23").Value = 110.375999
$ws.Range("I23").Value = 0.2820585476079611
$ws.Range("J23").Value = 0.2892552641013719
$ws.Range("O23").Value = 0.245697991654417
$ws.Range("P23").Value = 0.253077086664408
$ws.Range("Q23").Value = 2588.771619065883
$ws.Range("R23").Value = 23298.94457159295
$ws.Range("S23").Value = 0.06930121867623781
$ws.Range("T23").Value = 0.07320387954111911
$ws.Range("G24").Value = 36.79199966666667
$ws.Range("H24").Value = 110.375999
$ws.Range("I24").Value = 0.2820585476079611
$ws.Range("J24").Value = 0.2892552641013719
$ws.Range("M24").Value = 82.007665
$ws.Range("N24").Value = 246.022995
$ws.Range("O24").Value = 0.2863622109480123
$ws.Range("P24").Value = 0.2949625822722868
$ws.Range("Q24").Value = 3017.225983344112
$ws.Range("R24").Value = 27155.03385009701
$ws.Range("S24").Value = 0.08077090930980094
$ws.Range("T24").Value = 0.08531947963519294
$ws.Range("G25").Value = 36.79199966666667
$ws.Range("H25").Value = 110.375999
$ws.Range("I25").Value = 0.2820585476079611
$ws.Range("J25").Value = 0.2892552641013719
$ws.Range("M25").Value = 25.0501465
$ws.Range("N25").Value = 50.100293
$ws.Range("O25").Value = 0.0874724982879541
$ws.Range("P25").Value = 0.06006638442832619
$ws.Range("Q25").Value = 921.6449816779513
$ws.Range("R25").Value = 5529.869890067707
$ws.Range("S25").Value = 0.0246723658227402
$ws.Range("T25").Value = 0.01737451789143002
$ws.Range("G26").Value = 36.79199966666667
$ws.Range("H26").Value = 110.375999
$ws.Range("I26").Value = 0.2820585476079611
$ws.Range("J26").Value = 0.2892552641013719
$ws.Range("M26").Value = 95.05788666666668
$ws.Range("N26").Value = 285.17366
$ws.Range("O26").Value = 0.3319322235782747
$ws.Range("P26").Value = 0.3419012078510756
$ws.Range("Q26").Value = 3497.369734554038
$ws.Range("R26").Value = 31476.32761098635
$ws.Range("S26").Value = 0.09362432088676918
$ws.Range("T26").Value = 0.09889672417354092
